$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the consolidated dataset.
# It belongs chronologically right after the current row 36, so insert a
# fresh row at 37 — this pushes the former rows 37-71 down to 38-72 and
# grows the used range to A1:R72, exactly like Excel's own Insert Row.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new record's data.
$ws.Range("A37").Value = 11
$ws.Range("B37").Value = "Vega Monumental Concepción"
$ws.Range("C37").Value = "Bíobío"
$ws.Range("D37").Value = 44967
$ws.Range("E37").Value = 8
$ws.Range("F37").Value = 100112031
$ws.Range("G37").Value = "Poroto verde"
$ws.Range("H37").Value = "Magnum"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 100
$ws.Range("K37").Value = 20000
$ws.Range("L37").Value = 22000
$ws.Range("M37").Value = 21000
$ws.Range("N37").Value = "$/saco 25 kilos"
$ws.Range("O37").Value = "Región de O'Higgins"
$ws.Range("P37").Value = 840
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"
